$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in header C1 ---
$ws.Range("C1").Value = "Population"

# --- Add new header labels for columns L:Q ---
$ws.Range("L1").Value = "Unnamed: 0"
$ws.Range("M1").Value = "2023_WPR_ND_Risk_Exposure"
$ws.Range("N1").Value = "2023_WPR_ND_Risk_Vulnerability"
$ws.Range("O1").Value = "2023_WPR_ND_Risk_Susceptibility"
$ws.Range("P1").Value = "2023_WPR_ND_Lack_Adaptive_Capacities"
$ws.Range("Q1").Value = "2023_WPR_ND_Lack_Coping_Capacities"

# Match the header style (bold, bordered, centered) used by the existing A1:K1 headers
$ws.Range("K1").Copy()
$ws.Range("L1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply text after paste (PasteSpecial(formats) should not touch values, but set again to be safe)
$ws.Range("L1").Value = "Unnamed: 0"
$ws.Range("M1").Value = "2023_WPR_ND_Risk_Exposure"
$ws.Range("N1").Value = "2023_WPR_ND_Risk_Vulnerability"
$ws.Range("O1").Value = "2023_WPR_ND_Risk_Susceptibility"
$ws.Range("P1").Value = "2023_WPR_ND_Lack_Adaptive_Capacities"
$ws.Range("Q1").Value = "2023_WPR_ND_Lack_Coping_Capacities"

# --- Populate new climate-risk data columns L:Q for rows 2-43 ---
$data = New-Object 'object[,]' 42,6
$data[0,0]=0; $data[0,1]=2.29; $data[0,2]=16.95; $data[0,3]=11.51; $data[0,4]=37.51; $data[0,5]=11.28
$data[1,0]=1; $data[1,1]=0.15; $data[1,2]=8.84; $data[1,3]=4.44; $data[1,4]=19.53; $data[1,5]=7.98
$data[2,0]=2; $data[2,1]=0.05; $data[2,2]=11.14; $data[2,3]=6.49; $data[2,4]=36.57; $data[2,5]=5.83
$data[3,0]=3; $data[3,1]=1.84; $data[3,2]=18.49; $data[3,3]=8.029999999999999; $data[3,4]=28.91; $data[3,5]=27.21
$data[4,0]=4; $data[4,1]=0.34; $data[4,2]=20.28; $data[4,3]=14.13; $data[4,4]=51.26; $data[4,5]=11.51
$data[5,0]=5; $data[5,1]=0.3; $data[5,2]=19.15; $data[5,3]=17.48; $data[5,4]=46.46; $data[5,5]=8.65
$data[6,0]=6; $data[6,1]=1.57; $data[6,2]=14.48; $data[6,3]=8.539999999999999; $data[6,4]=37.03; $data[6,5]=9.6
$data[7,0]=7; $data[7,1]=1.02; $data[7,2]=12.43; $data[7,3]=7.26; $data[7,4]=35.48; $data[7,5]=7.45
$data[8,0]=8; $data[8,1]=$null; $data[8,2]=$null; $data[8,3]=$null; $data[8,4]=$null; $data[8,5]=$null
$data[9,0]=9; $data[9,1]=0.18; $data[9,2]=5.42; $data[9,3]=3.47; $data[9,4]=28.66; $data[9,5]=1.6
$data[10,0]=10; $data[10,1]=0.43; $data[10,2]=7.87; $data[10,3]=6.7; $data[10,4]=37.87; $data[10,5]=1.92
$data[11,0]=11; $data[11,1]=0.49; $data[11,2]=4.2; $data[11,3]=5.71; $data[11,4]=26.53; $data[11,5]=0.49
$data[12,0]=12; $data[12,1]=2.7; $data[12,2]=20.23; $data[12,3]=8.460000000000001; $data[12,4]=33.29; $data[12,5]=29.38
$data[13,0]=13; $data[13,1]=1.99; $data[13,2]=9.279999999999999; $data[13,3]=7.02; $data[13,4]=35.42; $data[13,5]=3.21
$data[14,0]=14; $data[14,1]=8.25; $data[14,2]=8.93; $data[14,3]=8.82; $data[14,4]=9.98; $data[14,5]=8.09
$data[15,0]=15; $data[15,1]=0.11; $data[15,2]=7.98; $data[15,3]=5.33; $data[15,4]=10.36; $data[15,5]=9.220000000000001
$data[16,0]=16; $data[16,1]=0.55; $data[16,2]=5.97; $data[16,3]=6.26; $data[16,4]=19.86; $data[16,5]=1.71
$data[17,0]=17; $data[17,1]=1.45; $data[17,2]=7.3; $data[17,3]=4.59; $data[17,4]=24.76; $data[17,5]=3.43
$data[18,0]=18; $data[18,1]=8.69; $data[18,2]=11.43; $data[18,3]=7.96; $data[18,4]=35.77; $data[18,5]=5.25
$data[19,0]=19; $data[19,1]=$null; $data[19,2]=$null; $data[19,3]=$null; $data[19,4]=$null; $data[19,5]=$null
$data[20,0]=20; $data[20,1]=$null; $data[20,2]=$null; $data[20,3]=$null; $data[20,4]=$null; $data[20,5]=$null
$data[21,0]=21; $data[21,1]=0.79; $data[21,2]=9.49; $data[21,3]=9.289999999999999; $data[21,4]=40.13; $data[21,5]=2.29
$data[22,0]=22; $data[22,1]=0.08; $data[22,2]=66.43000000000001; $data[22,3]=61.03; $data[22,4]=68.44; $data[22,5]=70.18000000000001
$data[23,0]=23; $data[23,1]=0.06; $data[23,2]=6.81; $data[23,3]=5.36; $data[23,4]=10.22; $data[23,5]=5.76
$data[24,0]=24; $data[24,1]=0.15; $data[24,2]=5.13; $data[24,3]=4.67; $data[24,4]=13.99; $data[24,5]=2.07
$data[25,0]=25; $data[25,1]=0.1; $data[25,2]=17.72; $data[25,3]=11.78; $data[25,4]=48.54; $data[25,5]=9.73
$data[26,0]=26; $data[26,1]=0.83; $data[26,2]=11.8; $data[26,3]=8.369999999999999; $data[26,4]=44.66; $data[26,5]=4.39
$data[27,0]=27; $data[27,1]=2.2; $data[27,2]=8.470000000000001; $data[27,3]=5.6; $data[27,4]=33.07; $data[27,5]=3.28
$data[28,0]=28; $data[28,1]=0.33; $data[28,2]=5.12; $data[28,3]=8.07; $data[28,4]=7.43; $data[28,5]=2.24
$data[29,0]=29; $data[29,1]=1.06; $data[29,2]=7.9; $data[29,3]=6.55; $data[29,4]=23.67; $data[29,5]=3.18
$data[30,0]=30; $data[30,1]=1.73; $data[30,2]=10.28; $data[30,3]=5.2; $data[30,4]=40.15; $data[30,5]=5.21
$data[31,0]=31; $data[31,1]=1.09; $data[31,2]=25.96; $data[31,3]=18.99; $data[31,4]=67.72; $data[31,5]=13.61
$data[32,0]=32; $data[32,1]=0.71; $data[32,2]=15.6; $data[32,3]=9.25; $data[32,4]=47.28; $data[32,5]=8.68
$data[33,0]=33; $data[33,1]=28.35; $data[33,2]=28.05; $data[33,3]=14.97; $data[33,4]=37.81; $data[33,5]=39
$data[34,0]=34; $data[34,1]=0.17; $data[34,2]=17.92; $data[34,3]=14.15; $data[34,4]=42.49; $data[34,5]=9.57
$data[35,0]=35; $data[35,1]=0.1; $data[35,2]=9.02; $data[35,3]=4.59; $data[35,4]=39.25; $data[35,5]=4.08
$data[36,0]=36; $data[36,1]=0.31; $data[36,2]=12.4; $data[36,3]=7.44; $data[36,4]=35.63; $data[36,5]=7.19
$data[37,0]=37; $data[37,1]=7.77; $data[37,2]=11.97; $data[37,3]=6.97; $data[37,4]=33.35; $data[37,5]=7.38
$data[38,0]=38; $data[38,1]=1.05; $data[38,2]=7.05; $data[38,3]=3.8; $data[38,4]=16.65; $data[38,5]=5.55
$data[39,0]=39; $data[39,1]=0.16; $data[39,2]=6.5; $data[39,3]=4.31; $data[39,4]=23.84; $data[39,5]=2.67
$data[40,0]=40; $data[40,1]=0.48; $data[40,2]=33.63; $data[40,3]=18.78; $data[40,4]=48.11; $data[40,5]=42.11
$data[41,0]=41; $data[41,1]=2.58; $data[41,2]=12.43; $data[41,3]=6.76; $data[41,4]=37.4; $data[41,5]=7.59

$ws.Range("L2:Q43").Value = $data
